$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '66.953.65'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '3.501.33'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.89'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.26'
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.585'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.132'
$ws.Range('E9').Value = '  +4.54%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.16'
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '4.103.84'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '29.46'
$ws.Range('E14').Value = '  +5.01%  '
$ws.Range('D15').Value = '66.953.63'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '3.495.63'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.28'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.26'
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '393.50'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.95'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.33'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.535'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.26'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.995'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.14'
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.06'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '23.67'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.38'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.62'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '164.19'
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.878'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.91'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.87'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '27.49'
$ws.Range('E39').Value = '  +3.73%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.65'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').Value = '2.848.16'
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0738'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '26.04'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.64'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.55'
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0301'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '338.85'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '34.79'
$ws.Range('E48').Value = '  +3.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.08'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.842'
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.43'
$ws.Range('E51').Value = '  -0.91%  '
